# Rename the last sheet ("Winter2016") to "2016-126" (CDSCC IF patch sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Winter2016")
$ws.Name = "2016-126"
